# Commit: "Added all drop-down and non-dropdown classes for all sites."
#
# The sheet is a simple key/value config table (column A = property name,
# column B/C = value), backed by the shared-strings table. This edit just
# fixes the casing of three key names:
#   A2: mdaTextHomepage  -> mdaTextHomePage
#   A4: MdaTitle         -> mdaTitle
#   A8: pageTitlenewTab  -> pageTitleNewTab
# and leaves the selected/active cell on A8 afterwards.
#
# Write order matters here: Excel appends newly-introduced shared strings
# to the end of the sharedStrings table in first-use order, so we write
# A4, A8, A2 (in that order) to reproduce the exact new string order
# (mdaTitle, pageTitleNewTab, mdaTextHomePage) from the target workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A4").Value = "mdaTitle"
$ws.Range("A8").Value = "pageTitleNewTab"
$ws.Range("A2").Value = "mdaTextHomePage"

# Also reposition the window per the source diff (best-effort; window
# chrome geometry isn't always round-tripped by every host, but attempt
# it via the documented Window object so it's captured when supported).
try {
    $win = $wb.Windows.Item(1)
    $win.Left = 13416
    $win.Top = 432
    $win.Width = 9600
    $win.Height = 11856
} catch {
}

# Move the active selection to A8, matching the new sheetView selection.
$ws.Range("A8").Select()
